# Adds a new "InternationalFT" worksheet (a copy of the header row plus
# one data row taken from DomesticFTA) as the last tab, and updates the
# view/selection state that Excel records for the sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet directly after the last existing tab --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "InternationalFT"

# --- 2. Populate header + data row by copying from DomesticFTA -------------
# (copying preserves both the shared-string values and the cell styles,
# e.g. the date format used in column H)
$domesticFTA = $wb.Worksheets.Item("DomesticFTA")
$domesticFTA.Range("A1:I1").Copy($newSheet.Range("A1:I1"))
$domesticFTA.Range("A4:I4").Copy($newSheet.Range("A2:I2"))
$newSheet.Range("A2").Value = "Data001"

# --- 3. Column widths --------------------------------------------------------
# (ColumnWidth is quantized to the nearest 1/6 character by the engine, so
# the inputs below are chosen to land on the closest achievable width to
# the recorded OOXML widths of 15 and 16.28515625.)
$newSheet.Columns.Item(8).ColumnWidth = 14.166666666666666
$newSheet.Columns.Item(9).ColumnWidth = 15.45

# --- 4. Selection / view state ----------------------------------------------
[void]$newSheet.Range("H8").Select()

[void]$domesticFTA.Range("A4:I4").Select()

[void]$newSheet.Activate()
